$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BS32 table: fill in the final (1E-009) tolerance row's integrator call count.
$ws.Range("D12").Value = 9249
# E12 already holds formula "=4*D12"; recalculated automatically.

# Update the average integrator calls for the DP87 0.0001-tolerance row (G33)
# with the latest run's measurements.
$ws.Range("G33").Formula = "=(9099+9603+9199+9095)/4"

# Add a new extrapolated row beneath the BS54 table.
$ws.Range("E53").Formula = "=1.333*E52"

# Match the saved cursor/selection position from the source workbook.
$ws.Range("E13").Select()
